# Update the parking register sheet: extend existing receipts with new
# exit times / amounts, and append new rows for longer parking-lot names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Recibo, Placa, Entrada, Salida, Vehiculo, Valor, Tiempo, Total
$rows = @(
    @{ R = 2;  A = "000076"; B = "AAA01"; C = "04/11/2024 17:57:00"; D = "05/11/2024 10:25:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 3;  A = "000077"; B = "BBB01"; C = "04/11/2024 17:58:00"; D = "05/11/2024 10:34:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 4;  A = "000078"; B = "CCC01"; C = "04/11/2024 17:58:00"; D = "05/11/2024 10:34:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 5;  A = "000079"; B = "DDD01"; C = "04/11/2024 17:58:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 6;  A = "000080"; B = "EEE01"; C = "04/11/2024 17:58:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 7;  A = "000081"; B = "FFF01"; C = "04/11/2024 17:59:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 8;  A = "000082"; B = "GGG01"; C = "04/11/2024 17:59:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 9;  A = "000083"; B = "HHH01"; C = "04/11/2024 17:59:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 },
    @{ R = 10; A = "000084"; B = "III01"; C = "04/11/2024 17:59:00"; D = "05/11/2024 10:35:00"; E = "Moto"; F = 7000; G = 0; H = 14000 }
)

foreach ($row in $rows) {
    $r = $row.R
    # Prefix the receipt number with an apostrophe so Excel keeps the
    # leading zeros and stores it as text instead of a number.
    $ws.Range("A$r").Value = "'" + $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}
